$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price & 1h volume change)
# Cells in column D whose new text looks like a plain number must be forced
# to Text format first, otherwise Excel auto-converts them to numeric values
# and the exact original text representation (leading/trailing zeros, etc.)
# would be lost.

$ws.Range("D2").Value = "29.922.97"
$ws.Range("E2").Value = "  +6.24%  "

$ws.Range("D3").Value = "1.877.07"
$ws.Range("E3").Value = "  +5.34%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.44"
$ws.Range("E5").Value = "  +1.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5003"
$ws.Range("E7").Value = "  +1.90%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.73"
$ws.Range("E8").Value = "  +7.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2857"
$ws.Range("E9").Value = "  +6.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06540"
$ws.Range("E10").Value = "  +4.22%  "

$ws.Range("D11").Value = "1.880.06"
$ws.Range("E11").Value = "  +5.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.09"
$ws.Range("E12").Value = "  +3.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07242"
$ws.Range("E13").Value = "  +3.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6663"
$ws.Range("E14").Value = "  +6.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.03"
$ws.Range("E15").Value = "  +6.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.826"
$ws.Range("E16").Value = "  +3.61%  "

$ws.Range("D17").Value = "29.940.63"
$ws.Range("E17").Value = "  +6.40%  "

$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("E19").Value = "  +6.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007517"
$ws.Range("E20").Value = "  +3.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").Value = "2.123.06"
$ws.Range("E22").Value = "  +5.53%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.767"
$ws.Range("E23").Value = "  +4.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.521"
$ws.Range("E24").Value = "  +5.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.016"
$ws.Range("E25").Value = "  +3.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.44"
$ws.Range("E26").Value = "  +3.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "134.86"
$ws.Range("E27").Value = "  +23.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.72"
$ws.Range("E28").Value = "  +6.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.950"
$ws.Range("E29").Value = "  +4.89%  "

$ws.Range("E30").Value = "  -0.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.195"
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08655"
$ws.Range("E32").Value = "  +4.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.894"
$ws.Range("E33").Value = "  +3.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05071"
$ws.Range("E34").Value = "  +3.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").Value = "  +5.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6918"
$ws.Range("E36").Value = "  +6.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  +2.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.297"
$ws.Range("E38").Value = "  +12.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.740"
$ws.Range("E39").Value = "  +5.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9600"
$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("E41").Value = "  +5.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.090"
$ws.Range("E42").Value = "  +3.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.62"
$ws.Range("E43").Value = "  +4.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.0000"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4210"
$ws.Range("E45").Value = "  +5.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.438"
$ws.Range("E46").Value = "  +3.66%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1255"
$ws.Range("E47").Value = "  +3.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05635"
$ws.Range("E48").Value = "  +3.63%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.38"
$ws.Range("E49").Value = "  +5.34%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.260"
$ws.Range("E50").Value = "  +3.19%  "

$ws.Range("E51").Value = "  +6.65%  "
